# Update "想去人数" (want-to-go count) figures for the two sheets that
# list exhibition events: "展览" (sheet1) and "全部类型" (sheet4).
# F2: 80 -> 82
# F3: 316 -> 318
# F4: 4443 -> 4475

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 82
    $ws.Range("F3").Value = 318
    $ws.Range("F4").Value = 4475
}
